$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "139 Highett St Apartment Complex Richmond",
    "3175 The Bays Aged Care Facility Hastings",
    "3563 Embracia Aged Care Reservoir",
    "Apartment Complex Fawkner",
    "Australian Lamb Colac East",
    "Bread Solutions Braeside",
    "CS Square Caroline Springs",
    "Carton Finishing Pty. Ltd. Campbellfield",
    "Cedar Meats Australia Brooklyn",
    "Community Kids Bayswater Early Education Centre Bayswater North",
    "Costco Wholesale Epping",
    "Ermha365 Residential Disability Care Services Doveton",
    "FedEx Station Melbourne Airport",
    "Green Leaves Early Learning Centre Highlands Craigieburn",
    "Guardian Childcare Caulfield",
    "Kool Kidz Childcare Narre Warren",
    "Lantmannen Unibake Australia Mordialloc",
    "Melbourne Assessment Prison West Melbourne",
    "MyCentre Childcare Broadmeadows",
    "Nido Early School Ascot Vale",
    "Nido Early School Glenroy",
    "Northern Health Northern Hospital Epping Emergency Department Tier 1B",
    "Northern Health The Northern Hospital Epping",
    "Social Gathering Warrnambool 28 September",
    "St Margaret's Primary School OSHC Maribyrnong",
    "St Vincents Hospital Emergency Department Melbourne",
    "The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B",
    "The Royal Talbot Rehabilitation Centre",
    "Visy Recycling Springvale",
    "Wallaby Childcare Wollert",
    "Werribee Mercy Hospital Emergency Department",
    "Western Health Footscray Hospital Emergency Department",
    "Western Health Sunshine Hospital Emergency Department"
)

$values = @(10, 14, 22, 10, 13, 13, 13, 12, 10, 18, 13, 10, 14, 14, 14, 16, 26, 10, 10, 29, 24, 52, 21, 17, 11, 35, 16, 12, 31, 16, 20, 10, 16)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
